# Updates cryptos list values (Price / Volume(1h) columns, plus a row-18/19
# coin swap) to match the latest scrape, per commit 'Updated cryptos list
# on Wed Mar 20 08:56:55 UTC 2024 with GitHub Actions'.
#
# All Price/Volume cells are stored as plain text in the sheet (no numeric
# cell type is used anywhere in this table), so every write below targets
# .Value directly. A handful of new Price strings (e.g. '525.32') look like
# plain decimals, and Excel would silently reinterpret them as numbers on
# assignment; for those cells we briefly force a text number format so the
# literal text is preserved, then restore the cell's original (default)
# style so no formatting residue is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '62.841.75'
$ws.Range('E2').Value = '  -0.96%  '
# Row 3
$ws.Range('D3').Value = '3.217.36'
$ws.Range('E3').Value = '  -1.05%  '
# Row 4
$ws.Range('E4').Value = '  -0.10%  '
# Row 5
$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '525.32'
$c.Style = "Normal"
$ws.Range('E5').Value = '  +2.90%  '
# Row 6
$c = $ws.Range('D6')
$c.NumberFormat = "@"
$c.Value = '171.94'
$c.Style = "Normal"
$ws.Range('E6').Value = '  -2.51%  '
# Row 7
$c = $ws.Range('D7')
$c.NumberFormat = "@"
$c.Value = '0.593'
$c.Style = "Normal"
$ws.Range('E7').Value = '  +2.07%  '
# Row 8
$ws.Range('E8').Value = '  -0.06%  '
# Row 9
$ws.Range('D9').Value = '3.217.98'
$ws.Range('E9').Value = '  -0.88%  '
# Row 10
$c = $ws.Range('D10')
$c.NumberFormat = "@"
$c.Value = '0.604'
$c.Style = "Normal"
$ws.Range('E10').Value = '  +0.48%  '
# Row 11
$c = $ws.Range('D11')
$c.NumberFormat = "@"
$c.Value = '52.94'
$c.Style = "Normal"
$ws.Range('E11').Value = '  -6.22%  '
# Row 12
$ws.Range('E12').Value = '  +4.13%  '
# Row 13
$ws.Range('E13').Value = '  +2.11%  '
# Row 14
$ws.Range('E14').Value = '  +2.61%  '
# Row 15
$ws.Range('D15').Value = '3.739.63'
$ws.Range('E15').Value = '  -1.74%  '
# Row 16
$ws.Range('E16').Value = '  -4.24%  '
# Row 17
$ws.Range('D17').Value = '3.217.36'
$ws.Range('E17').Value = '  -1.85%  '
# Row 18
$ws.Range('B18').Value = 'WrappedBTC'
$ws.Range('C18').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D18').Value = '62.811.69'
$ws.Range('E18').Value = '  -0.82%  '
# Row 19
$ws.Range('B19').Value = 'Chainlink'
$ws.Range('C19').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$c = $ws.Range('D19')
$c.NumberFormat = "@"
$c.Value = '17.17'
$c.Style = "Normal"
$ws.Range('E19').Value = '  +2.24%  '
# Row 20
$c = $ws.Range('D20')
$c.NumberFormat = "@"
$c.Value = '11.01'
$c.Style = "Normal"
$ws.Range('E20').Value = '  +4.23%  '
# Row 21
$ws.Range('E21').Value = '  +4.43%  '
# Row 22
$c = $ws.Range('D22')
$c.NumberFormat = "@"
$c.Value = '364.60'
$c.Style = "Normal"
$ws.Range('E22').Value = '  +0.59%  '
# Row 23
$c = $ws.Range('D23')
$c.NumberFormat = "@"
$c.Value = '3.76'
$c.Style = "Normal"
$ws.Range('E23').Value = '  +5.00%  '
# Row 24
$c = $ws.Range('D24')
$c.NumberFormat = "@"
$c.Value = '81.08'
$c.Style = "Normal"
$ws.Range('E24').Value = '  +3.60%  '
# Row 25
$ws.Range('E25').Value = '  +4.11%  '
# Row 26
$c = $ws.Range('D26')
$c.NumberFormat = "@"
$c.Value = '3.89'
$c.Style = "Normal"
$ws.Range('E26').Value = '  +5.60%  '
# Row 27
$ws.Range('E27').Value = '  +3.03%  '
# Row 28
$c = $ws.Range('D28')
$c.NumberFormat = "@"
$c.Value = '2.63'
$c.Style = "Normal"
$ws.Range('E28').Value = '  +2.56%  '
# Row 29
$c = $ws.Range('D29')
$c.NumberFormat = "@"
$c.Value = '11.24'
$c.Style = "Normal"
$ws.Range('E29').Value = '  +3.42%  '
# Row 30
$c = $ws.Range('D30')
$c.NumberFormat = "@"
$c.Value = '8.11'
$c.Style = "Normal"
$ws.Range('E30').Value = '  +0.66%  '
# Row 31
$c = $ws.Range('D31')
$c.NumberFormat = "@"
$c.Value = '28.41'
$c.Style = "Normal"
$ws.Range('E31').Value = '  +2.31%  '
# Row 32
$c = $ws.Range('D32')
$c.NumberFormat = "@"
$c.Value = '628.96'
$c.Style = "Normal"
$ws.Range('E32').Value = '  -1.62%  '
# Row 33
$c = $ws.Range('D33')
$c.NumberFormat = "@"
$c.Value = '6.42'
$c.Style = "Normal"
$ws.Range('E33').Value = '  -0.54%  '
# Row 34
$c = $ws.Range('D34')
$c.NumberFormat = "@"
$c.Value = '11.21'
$c.Style = "Normal"
$ws.Range('E34').Value = '  +3.78%  '
# Row 35
$ws.Range('E35').Value = '  +5.40%  '
# Row 36
$c = $ws.Range('D36')
$c.NumberFormat = "@"
$c.Value = '56.81'
$c.Style = "Normal"
$ws.Range('E36').Value = '  -1.97%  '
# Row 38
$c = $ws.Range('D38')
$c.NumberFormat = "@"
$c.Value = '36.63'
$c.Style = "Normal"
$ws.Range('E38').Value = '  +4.53%  '
# Row 39
$ws.Range('E39').Value = '  +3.14%  '
# Row 40
$c = $ws.Range('D40')
$c.NumberFormat = "@"
$c.Value = '0.997'
$c.Style = "Normal"
$ws.Range('E40').Value = '  -0.21%  '
# Row 41
$ws.Range('D41').Value = '0.0₃0709'
$ws.Range('E41').Value = '  +16.61%  '
# Row 42
$ws.Range('E42').Value = '  +2.42%  '
# Row 43
$ws.Range('D43').Value = '2.868.12'
$ws.Range('E43').Value = '  +3.80%  '
# Row 44
$c = $ws.Range('D44')
$c.NumberFormat = "@"
$c.Value = '2.54'
$c.Style = "Normal"
$ws.Range('E44').Value = '  +14.22%  '
# Row 45
$c = $ws.Range('D45')
$c.NumberFormat = "@"
$c.Value = '2.69'
$c.Style = "Normal"
$ws.Range('E45').Value = '  +4.56%  '
# Row 46
$c = $ws.Range('D46')
$c.NumberFormat = "@"
$c.Value = '2.91'
$c.Style = "Normal"
$ws.Range('E46').Value = '  +14.20%  '
# Row 47
$ws.Range('E47').Value = '  +4.94%  '
# Row 48
$ws.Range('E48').Value = '  -1.60%  '
# Row 49
$c = $ws.Range('D49')
$c.NumberFormat = "@"
$c.Value = '3.02'
$c.Style = "Normal"
$ws.Range('E49').Value = '  +9.56%  '
# Row 50
$ws.Range('E50').Value = '  +2.89%  '
# Row 51
$c = $ws.Range('D51')
$c.NumberFormat = "@"
$c.Value = '135.28'
$c.Style = "Normal"
$ws.Range('E51').Value = '  +2.12%  '
